# From v1.0.3 to v1.1
# Swap the "Steps" / "Expected Results" content between the TC3 block
# (rows 21-25) and the TC4 block (rows 28-32), so that the "cancelar
# diária" step now belongs to TC3 and the "detalhar diária" step now
# belongs to TC4 (the TC3 / TC4 labels themselves stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Steps = $ws.Range("B25").Value2
$tc3Expected = $ws.Range("D25").Value2
$tc4Steps = $ws.Range("B32").Value2
$tc4Expected = $ws.Range("D32").Value2

$ws.Range("B25").Value2 = $tc4Steps
$ws.Range("D25").Value2 = $tc4Expected
$ws.Range("B32").Value2 = $tc3Steps
$ws.Range("D32").Value2 = $tc3Expected
